$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row Right count
$ws.Range("B11").Value = 5

# Update "Total" row Right count and Corr/total label
$ws.Range("B12").Value = 130
$ws.Range("E12").Value = "130/140"
